$wb = $excel.ActiveWorkbook

# ---- Sheet: LP1912 ----
$ws = $wb.Worksheets.Item('LP1912')
$ws.Range('A2').Value = 'Última actualización: 11:00:36'
$ws.Range('A3').Value = 'Total filas: 171'
$ws.Range('A24').Value = '03:42:43'
$ws.Range('C24').Value = '14_ABASTO'
$ws.Range('D24').Value = 113
$ws.Range('A25').Value = '04:17:03'
$ws.Range('C25').Value = '215B_EL PATO'
$ws.Range('D25').Value = 78
$ws.Range('C55').Value = '84_COLONIA URQUIZA-ESC 49'
$ws.Range('C57').Value = '11_ETCHEVERRY'
$ws.Range('A85').Value = '07:50:23'
$ws.Range('C85').Value = '10_OLMOS'
$ws.Range('D85').Value = 63
$ws.Range('A86').Value = '07:17:57'
$ws.Range('C86').Value = '17_ROMERO'
$ws.Range('D86').Value = 96
$ws.Range('C87').Value = '225_HARAS DEL SUR'
$ws.Range('C88').Value = '17_ROMERO'
$ws.Range('C98').Value = '23_HERNANDEZ'
$ws.Range('C99').Value = '26_HERNANDEZ'
$ws.Range('A100').Value = '08:52:26'
$ws.Range('C100').Value = '16_SANTA ANA'
$ws.Range('D100').Value = 30
$ws.Range('A101').Value = '08:39:38'
$ws.Range('C101').Value = '17_ROMERO'
$ws.Range('D101').Value = 43
$ws.Range('A102').Value = '07:50:23'
$ws.Range('C102').Value = '17_ROMERO'
$ws.Range('D102').Value = 93
$ws.Range('A103').Value = '08:52:26'
$ws.Range('C103').Value = '11_ETCHEVERRY'
$ws.Range('D103').Value = 31
$ws.Range('A112').Value = '08:39:38'
$ws.Range('C112').Value = '10_OLMOS'
$ws.Range('D112').Value = 63
$ws.Range('A113').Value = '09:28:24'
$ws.Range('C113').Value = '215C_EL PATO'
$ws.Range('D113').Value = 14
$ws.Range('C128').Value = '16_SANTA ANA'
$ws.Range('C129').Value = '23_HERNANDEZ'
$ws.Range('A138').Value = '11:00:36'
$ws.Range('C138').Value = '16_SANTA ANA'
$ws.Range('D138').Value = 1
$ws.Range('A139').Value = '11:00:36'
$ws.Range('B139').Value = '11:01'
$ws.Range('D139').Value = 1
$ws.Range('A140').Value = '09:28:24'
$ws.Range('B140').Value = '11:02'
$ws.Range('C140').Value = '215C_EL PATO'
$ws.Range('D140').Value = 94
$ws.Range('A141').Value = '11:00:36'
$ws.Range('B141').Value = '11:02'
$ws.Range('C141').Value = '11_ETCHEVERRY'
$ws.Range('D141').Value = 2
$ws.Range('B142').Value = '11:03'
$ws.Range('C142').Value = '11_ETCHEVERRY'
$ws.Range('D142').Value = 38
$ws.Range('A143').Value = '11:00:36'
$ws.Range('B143').Value = '11:04'
$ws.Range('C143').Value = '23_HERNANDEZ'
$ws.Range('D143').Value = 4
$ws.Range('A144').Value = '11:00:36'
$ws.Range('B144').Value = '11:06'
$ws.Range('C144').Value = '16_P MOR-167 Y 521'
$ws.Range('D144').Value = 6
$ws.Range('A145').Value = '11:00:36'
$ws.Range('B145').Value = '11:11'
$ws.Range('C145').Value = '10_OLMOS'
$ws.Range('D145').Value = 11
$ws.Range('A146').Value = '11:00:36'
$ws.Range('B146').Value = '11:12'
$ws.Range('C146').Value = '15_ABASTO'
$ws.Range('D146').Value = 12
$ws.Range('A147').Value = '11:00:36'
$ws.Range('B147').Value = '11:16'
$ws.Range('C147').Value = '16_SANTA ANA'
$ws.Range('D147').Value = 16
$ws.Range('A148').Value = '11:00:36'
$ws.Range('B148').Value = '11:19'
$ws.Range('C148').Value = '86_EST CHICA-ESC AGRARIA'
$ws.Range('D148').Value = 19
$ws.Range('B149').Value = '11:20'
$ws.Range('C149').Value = '26_HERNANDEZ'
$ws.Range('D149').Value = 55
$ws.Range('A150').Value = '11:00:36'
$ws.Range('B150').Value = '11:21'
$ws.Range('C150').Value = '26_HERNANDEZ'
$ws.Range('D150').Value = 21
$ws.Range('A151').Value = '11:00:36'
$ws.Range('B151').Value = '11:22'
$ws.Range('C151').Value = '17_ROMERO'
$ws.Range('D151').Value = 22
$ws.Range('B152').Value = '11:26'
$ws.Range('C152').Value = '225_C ROCA-H SUR'
$ws.Range('D152').Value = 61
$ws.Range('A153').Value = '11:00:36'
$ws.Range('B153').Value = '11:27'
$ws.Range('C153').Value = '225_C ROCA-H SUR'
$ws.Range('D153').Value = 27
$ws.Range('A154').Value = '11:00:36'
$ws.Range('B154').Value = '11:32'
$ws.Range('C154').Value = '81_EL PELIGRO'
$ws.Range('D154').Value = 32
$ws.Range('A155').Value = '11:00:36'
$ws.Range('B155').Value = '11:34'
$ws.Range('C155').Value = '23_HERNANDEZ'
$ws.Range('D155').Value = 34
$ws.Range('A156').Value = '11:00:36'
$ws.Range('B156').Value = '11:35'
$ws.Range('C156').Value = '11_ETCHEVERRY'
$ws.Range('D156').Value = 35
$ws.Range('B157').Value = '11:41'
$ws.Range('C157').Value = '17_ROMERO'
$ws.Range('D157').Value = 76
$ws.Range('A158').Value = '11:00:36'
$ws.Range('B158').Value = '11:42'
$ws.Range('C158').Value = '17_ROMERO'
$ws.Range('D158').Value = 42
$ws.Range('A159').Value = '11:00:36'
$ws.Range('B159').Value = '11:43'
$ws.Range('C159').Value = '10_OLMOS'
$ws.Range('D159').Value = 43
$ws.Range('A160').Value = '10:25:56'
$ws.Range('B160').Value = '11:49'
$ws.Range('C160').Value = '15_ABASTO'
$ws.Range('D160').Value = 84
$ws.Range('E160').Value = 'LP1912'
$ws.Range('A161').Value = '11:00:36'
$ws.Range('B161').Value = '11:51'
$ws.Range('C161').Value = '215B_EL PATO'
$ws.Range('D161').Value = 51
$ws.Range('E161').Value = 'LP1912'
$ws.Range('A162').Value = '11:00:36'
$ws.Range('B162').Value = '11:52'
$ws.Range('C162').Value = '15_ABASTO'
$ws.Range('D162').Value = 52
$ws.Range('E162').Value = 'LP1912'
$ws.Range('A163').Value = '10:25:56'
$ws.Range('B163').Value = '11:58'
$ws.Range('C163').Value = '225_GOMEZ'
$ws.Range('D163').Value = 93
$ws.Range('E163').Value = 'LP1912'
$ws.Range('A164').Value = '11:00:36'
$ws.Range('B164').Value = '11:59'
$ws.Range('C164').Value = '225_GOMEZ'
$ws.Range('D164').Value = 59
$ws.Range('E164').Value = 'LP1912'
$ws.Range('A165').Value = '11:00:36'
$ws.Range('B165').Value = '12:02'
$ws.Range('C165').Value = '84_COLONIA URQUIZA-ESC 49'
$ws.Range('D165').Value = 62
$ws.Range('E165').Value = 'LP1912'
$ws.Range('A166').Value = '11:00:36'
$ws.Range('B166').Value = '12:06'
$ws.Range('C166').Value = '16_P MOR-SANTA ANA'
$ws.Range('D166').Value = 66
$ws.Range('E166').Value = 'LP1912'
$ws.Range('A167').Value = '11:00:36'
$ws.Range('B167').Value = '12:06'
$ws.Range('C167').Value = '14_ABASTO'
$ws.Range('D167').Value = 66
$ws.Range('E167').Value = 'LP1912'
$ws.Range('A168').Value = '11:00:36'
$ws.Range('B168').Value = '12:13'
$ws.Range('C168').Value = '10_OLMOS'
$ws.Range('D168').Value = 73
$ws.Range('E168').Value = 'LP1912'
$ws.Range('A169').Value = '11:00:36'
$ws.Range('B169').Value = '12:20'
$ws.Range('C169').Value = '215A_EL PATO'
$ws.Range('D169').Value = 80
$ws.Range('E169').Value = 'LP1912'
$ws.Range('A170').Value = '10:25:56'
$ws.Range('B170').Value = '12:20'
$ws.Range('C170').Value = '26_HERNANDEZ'
$ws.Range('D170').Value = 115
$ws.Range('E170').Value = 'LP1912'
$ws.Range('A171').Value = '10:25:56'
$ws.Range('B171').Value = '12:20'
$ws.Range('C171').Value = '14_ABASTO'
$ws.Range('D171').Value = 115
$ws.Range('E171').Value = 'LP1912'
$ws.Range('A172').Value = '11:00:36'
$ws.Range('B172').Value = '12:21'
$ws.Range('C172').Value = '26_HERNANDEZ'
$ws.Range('D172').Value = 81
$ws.Range('E172').Value = 'LP1912'
$ws.Range('A173').Value = '11:00:36'
$ws.Range('B173').Value = '12:21'
$ws.Range('C173').Value = '14_ABASTO'
$ws.Range('D173').Value = 81
$ws.Range('E173').Value = 'LP1912'
$ws.Range('A174').Value = '11:00:36'
$ws.Range('B174').Value = '12:36'
$ws.Range('C174').Value = '27_EL RETIRO'
$ws.Range('D174').Value = 96
$ws.Range('E174').Value = 'LP1912'
$ws.Range('A175').Value = '11:00:36'
$ws.Range('B175').Value = '12:38'
$ws.Range('C175').Value = '17_179 Y 38'
$ws.Range('D175').Value = 98
$ws.Range('E175').Value = 'LP1912'
$ws.Range('A176').Value = '11:00:36'
$ws.Range('B176').Value = '12:48'
$ws.Range('C176').Value = '11_ETCHEVERRY'
$ws.Range('D176').Value = 108
$ws.Range('E176').Value = 'LP1912'

# ---- Sheet: LP1912-215 ----
$ws = $wb.Worksheets.Item('LP1912-215')
$ws.Range('A2').Value = 'Última actualización: 11:00:36'
$ws.Range('A28').Value = '11:00:36'
$ws.Range('D28').Value = 1
$ws.Range('A30').Value = '11:00:36'
$ws.Range('D30').Value = 51
$ws.Range('A31').Value = '11:00:36'
$ws.Range('D31').Value = 80

# ---- Sheet: 6203-6173 ----
$ws = $wb.Worksheets.Item('6203-6173')
$ws.Range('A2').Value = 'Última actualización: 11:00:36'
$ws.Range('A3').Value = 'Total filas: 32'
$ws.Range('A33').Value = '11:00:36'
$ws.Range('D33').Value = 13
$ws.Range('A36').Value = '11:00:36'
$ws.Range('B36').Value = '12:04'
$ws.Range('C36').Value = '215A_LA PLATA'
$ws.Range('D36').Value = 64
$ws.Range('E36').Value = 'L6173'
$ws.Range('A37').Value = '11:00:36'
$ws.Range('B37').Value = '12:53'
$ws.Range('C37').Value = '215C_LA PLATA'
$ws.Range('D37').Value = 113
$ws.Range('E37').Value = 'L6203'

